$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new entry was logged for the "September" group which pushes all the
# existing rows (from row 29 downward) down by one row. Insert a new row
# at row 29 to reproduce that shift.
$ws.Rows("29:29").Insert()

# The other columns of the freshly inserted row are blank, matching the
# rest of the sheet's rows (every cell is present, just empty). Touching
# the (no-op) Style property materializes the cell without altering its
# appearance.
for ($col = 1; $col -le 25; $col++) {
    $ws.Cells.Item(29, $col).Style = "Normal"
}

# Populate the newly inserted row with the new September entry.
$ws.Range("R29").Value = "internet"
$ws.Range("S29").Value = "2024-09-03 19:58:18"
